$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.388.50"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "2.301.27"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "317.15"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").Value = "104.03"
$ws.Range("E6").Value = "  -2.32%  "

$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  -0.53%  "

$ws.Range("D11").Value = "0.0909"

$ws.Range("D12").Value = "8.36"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("D14").Value = "0.968"
$ws.Range("E14").Value = "  -1.39%  "

$ws.Range("D15").Value = "15.32"
$ws.Range("E15").Value = "  -1.76%  "

$ws.Range("D16").Value = "2.649.01"
$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("D17").Value = "2.309.05"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "42.366.56"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").Value = "7.49"
$ws.Range("E19").Value = "  -3.15%  "

$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").Value = "73.36"
$ws.Range("E21").Value = "  -1.81%  "

$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").Value = "3.54"
$ws.Range("E22").Value = "  +1.34%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "276.37"
$ws.Range("E23").Value = "  +5.98%  "

$ws.Range("D24").Value = "10.94"
$ws.Range("E24").Value = "  +16.97%  "

$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  -1.21%  "

$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  -1.59%  "

$ws.Range("D28").Value = "2.40"
$ws.Range("E28").Value = "  +5.45%  "

$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").Value = "35.84"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("D31").Value = "165.01"
$ws.Range("E31").Value = "  +0.87%  "

$ws.Range("D32").Value = "0.0873"
$ws.Range("E32").Value = "  -3.51%  "

$ws.Range("D33").Value = "5.88"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("E34").Value = "  +4.36%  "

$ws.Range("E35").Value = "  -10.73%  "

$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").Value = "0.0370"
$ws.Range("E37").Value = "  +4.80%  "

$ws.Range("E38").Value = "  +1.50%  "

$ws.Range("E39").Value = "  +2.89%  "

$ws.Range("D40").Value = "2.76"
$ws.Range("E40").Value = "  -0.98%  "

$ws.Range("E41").Value = "  +1.74%  "

$ws.Range("D42").Value = "70.05"
$ws.Range("E42").Value = "  -2.24%  "

$ws.Range("D43").Value = "0.228"
$ws.Range("E43").Value = "  -0.90%  "

$ws.Range("D45").Value = "93.65"
$ws.Range("E45").Value = "  -3.88%  "

$ws.Range("D46").Value = "82.05"
$ws.Range("E46").Value = "  +10.34%  "

$ws.Range("E47").Value = "  -2.87%  "

$ws.Range("D48").Value = "112.70"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").Value = "8.91"
$ws.Range("E49").Value = "  -1.38%  "

$ws.Range("D50").Value = "1.601.04"
$ws.Range("E50").Value = "  +2.90%  "

$ws.Range("D51").Value = "5.12"
$ws.Range("E51").Value = "  -4.05%  "
